{"js": "/*\n * Applies the per-paragraph text replacements described by the diff:\n * the title date line and each of the 100 table-cell arithmetic\n * expressions (in document order: title paragraph first, then the\n * table's rows left-to-right, top-to-bottom) get their text swapped\n * for the new value, while leaving every other property (fonts,\n * sizes, alignment, table layout, etc.) untouched.\n */\nconst REPLACEMENTS = [\n  [\"2025-07-26 Saturday\", \"2025-07-27 Sunday\"],\n  [\"46+35=81\", \"44+17=61\"],\n  [\"20-15=5\", \"51-18=33\"],\n  [\"59+7=66\", \"72-46=26\"],\n  [\"63-49=14\", \"6+78=84\"],\n  [\"71-17=54\", \"29+33=62\"],\n  [\"66+16=82\", \"16+29=45\"],\n  [\"72-15=57\", \"34+59=93\"],\n  [\"25+16=41\", \"39+6=45\"],\n  [\"30-13=17\", \"79+14=93\"],\n  [\"91-43=48\", \"40-36=4\"],\n  [\"95-29=66\", \"30-12=18\"],\n  [\"46-28=18\", \"74-45=29\"],\n  [\"39+35=74\", \"30-29=1\"],\n  [\"70-14=56\", \"45+27=72\"],\n  [\"18+76=94\", \"95-78=17\"],\n  [\"36+49=85\", \"83-35=48\"],\n  [\"91-13=78\", \"71-3=68\"],\n  [\"34-25=9\", \"28+29=57\"],\n  [\"32-14=18\", \"80-64=16\"],\n  [\"94-5=89\", \"64-17=47\"],\n  [\"61-57=4\", \"51-29=22\"],\n  [\"86-69=17\", \"25+56=81\"],\n  [\"9+38=47\", \"86+6=92\"],\n  [\"65-19=46\", \"91-35=56\"],\n  [\"60-47=13\", \"49+46=95\"],\n  [\"82-73=9\", \"76-37=39\"],\n  [\"59+18=77\", \"68+16=84\"],\n  [\"58+38=96\", \"28+34=62\"],\n  [\"65+9=74\", \"71-38=33\"],\n  [\"40-26=14\", \"7+9=16\"],\n  [\"81-24=57\", \"93-19=74\"],\n  [\"46-19=27\", \"84-37=47\"],\n  [\"27+34=61\", \"39+4=43\"],\n  [\"73+8=81\", \"26+58=84\"],\n  [\"90-42=48\", \"31-15=16\"],\n  [\"61-48=13\", \"27+66=93\"],\n  [\"81-26=55\", \"16+66=82\"],\n  [\"36+16=52\", \"90-36=54\"],\n  [\"17+68=85\", \"73-54=19\"],\n  [\"53+39=92\", \"70-33=37\"],\n  [\"39+49=88\", \"78-59=19\"],\n  [\"23+48=71\", \"29+69=98\"],\n  [\"24+8=32\", \"70-42=28\"],\n  [\"8+36=44\", \"26+45=71\"],\n  [\"39+15=54\", \"94-26=68\"],\n  [\"64-16=48\", \"42-13=29\"],\n  [\"63-56=7\", \"27+39=66\"],\n  [\"31-28=3\", \"82-37=45\"],\n  [\"27+5=32\", \"78-59=19\"],\n  [\"47+25=72\", \"92-77=15\"],\n  [\"29+42=71\", \"18+49=67\"],\n  [\"82-45=37\", \"6+48=54\"],\n  [\"47-19=28\", \"66+6=72\"],\n  [\"4+49=53\", \"19+6=25\"],\n  [\"45+17=62\", \"19+15=34\"],\n  [\"26+46=72\", \"85-9=76\"],\n  [\"73-18=55\", \"40-31=9\"],\n  [\"37-18=19\", \"82-34=48\"],\n  [\"95-59=36\", \"28+5=33\"],\n  [\"18+57=75\", \"9+53=62\"],\n  [\"97-59=38\", \"6+79=85\"],\n  [\"57+25=82\", \"90-82=8\"],\n  [\"64-58=6\", \"65+8=73\"],\n  [\"49+17=66\", \"91-62=29\"],\n  [\"67+6=73\", \"74-35=39\"],\n  [\"6+6=12\", \"48+18=66\"],\n  [\"9+12=21\", \"16+58=74\"],\n  [\"77+8=85\", \"58-49=9\"],\n  [\"47+48=95\", \"83-78=5\"],\n  [\"84-25=59\", \"90-69=21\"],\n  [\"36+59=95\", \"71-52=19\"],\n  [\"20-13=7\", \"59+25=84\"],\n  [\"19+43=62\", \"18+6=24\"],\n  [\"94-86=8\", \"27+69=96\"],\n  [\"42+9=51\", \"9+29=38\"],\n  [\"8+79=87\", \"8+68=76\"],\n  [\"21-12=9\", \"28+55=83\"],\n  [\"54-15=39\", \"85-68=17\"],\n  [\"45-36=9\", \"16+35=51\"],\n  [\"66+16=82\", \"8+46=54\"],\n  [\"33-29=4\", \"67+7=74\"],\n  [\"18+55=73\", \"32-26=6\"],\n  [\"80-77=3\", \"7+74=81\"],\n  [\"87+5=92\", \"26+5=31\"],\n  [\"4+59=63\", \"80-7=73\"],\n  [\"60-53=7\", \"16+29=45\"],\n  [\"65-37=28\", \"45+29=74\"],\n  [\"54-35=19\", \"38+9=47\"],\n  [\"2+39=41\", \"26+59=85\"],\n  [\"81-78=3\", \"30-27=3\"],\n  [\"91-65=26\", \"46+16=62\"],\n  [\"17+68=85\", \"27+39=66\"],\n  [\"13+69=82\", \"74+17=91\"],\n  [\"33+9=42\", \"35+59=94\"],\n  [\"85-17=68\", \"24+48=72\"],\n  [\"27+28=55\", \"31-2=29\"],\n  [\"76-48=28\", \"81-72=9\"],\n  [\"90-83=7\", \"15+26=41\"],\n  [\"27+58=85\", \"24-15=9\"],\n  [\"85-46=39\", \"85+8=93\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    \"Expected \" + REPLACEMENTS.length + \" paragraphs, found \" + items.length\n  );\n}\n\n// Read current text for every paragraph first so we can verify we are\n// editing the expected paragraph before mutating it.\nitems.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nfor (let i = 0; i < items.length; i++) {\n  const oldText = REPLACEMENTS[i][0];\n  const newText = REPLACEMENTS[i][1];\n  const paragraph = items[i];\n  if (paragraph.text !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \": expected '\" + oldText + \"' but found '\" + paragraph.text + \"'\"\n    );\n  }\n  if (oldText === newText) {\n    continue;\n  }\n  // getRange() + insertText(..., \"Replace\") swaps the run's text while\n  // preserving the existing run/paragraph formatting (fonts, size,\n  // justification, etc.)\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Applies the per-paragraph text replacements described by the diff:\n# the title date line (index 0) and each of the 100 table-cell\n# arithmetic expressions (indices 1..100, in row-major reading order:\n# row 1 col 1, row 1 col 2, ... row 20 col 5) get their text swapped\n# for the new value, while leaving every other property (fonts, sizes,\n# alignment, table layout, etc.) untouched.\n$Replacements = @(\n    @('2025-07-26 Saturday', '2025-07-27 Sunday'),\n    @('46+35=81', '44+17=61'),\n    @('20-15=5', '51-18=33'),\n    @('59+7=66', '72-46=26'),\n    @('63-49=14', '6+78=84'),\n    @('71-17=54', '29+33=62'),\n    @('66+16=82', '16+29=45'),\n    @('72-15=57', '34+59=93'),\n    @('25+16=41', '39+6=45'),\n    @('30-13=17', '79+14=93'),\n    @('91-43=48', '40-36=4'),\n    @('95-29=66', '30-12=18'),\n    @('46-28=18', '74-45=29'),\n    @('39+35=74', '30-29=1'),\n    @('70-14=56', '45+27=72'),\n    @('18+76=94', '95-78=17'),\n    @('36+49=85', '83-35=48'),\n    @('91-13=78', '71-3=68'),\n    @('34-25=9', '28+29=57'),\n    @('32-14=18', '80-64=16'),\n    @('94-5=89', '64-17=47'),\n    @('61-57=4', '51-29=22'),\n    @('86-69=17', '25+56=81'),\n    @('9+38=47', '86+6=92'),\n    @('65-19=46', '91-35=56'),\n    @('60-47=13', '49+46=95'),\n    @('82-73=9', '76-37=39'),\n    @('59+18=77', '68+16=84'),\n    @('58+38=96', '28+34=62'),\n    @('65+9=74', '71-38=33'),\n    @('40-26=14', '7+9=16'),\n    @('81-24=57', '93-19=74'),\n    @('46-19=27', '84-37=47'),\n    @('27+34=61', '39+4=43'),\n    @('73+8=81', '26+58=84'),\n    @('90-42=48', '31-15=16'),\n    @('61-48=13', '27+66=93'),\n    @('81-26=55', '16+66=82'),\n    @('36+16=52', '90-36=54'),\n    @('17+68=85', '73-54=19'),\n    @('53+39=92', '70-33=37'),\n    @('39+49=88', '78-59=19'),\n    @('23+48=71', '29+69=98'),\n    @('24+8=32', '70-42=28'),\n    @('8+36=44', '26+45=71'),\n    @('39+15=54', '94-26=68'),\n    @('64-16=48', '42-13=29'),\n    @('63-56=7', '27+39=66'),\n    @('31-28=3', '82-37=45'),\n    @('27+5=32', '78-59=19'),\n    @('47+25=72', '92-77=15'),\n    @('29+42=71', '18+49=67'),\n    @('82-45=37', '6+48=54'),\n    @('47-19=28', '66+6=72'),\n    @('4+49=53', '19+6=25'),\n    @('45+17=62', '19+15=34'),\n    @('26+46=72', '85-9=76'),\n    @('73-18=55', '40-31=9'),\n    @('37-18=19', '82-34=48'),\n    @('95-59=36', '28+5=33'),\n    @('18+57=75', '9+53=62'),\n    @('97-59=38', '6+79=85'),\n    @('57+25=82', '90-82=8'),\n    @('64-58=6', '65+8=73'),\n    @('49+17=66', '91-62=29'),\n    @('67+6=73', '74-35=39'),\n    @('6+6=12', '48+18=66'),\n    @('9+12=21', '16+58=74'),\n    @('77+8=85', '58-49=9'),\n    @('47+48=95', '83-78=5'),\n    @('84-25=59', '90-69=21'),\n    @('36+59=95', '71-52=19'),\n    @('20-13=7', '59+25=84'),\n    @('19+43=62', '18+6=24'),\n    @('94-86=8', '27+69=96'),\n    @('42+9=51', '9+29=38'),\n    @('8+79=87', '8+68=76'),\n    @('21-12=9', '28+55=83'),\n    @('54-15=39', '85-68=17'),\n    @('45-36=9', '16+35=51'),\n    @('66+16=82', '8+46=54'),\n    @('33-29=4', '67+7=74'),\n    @('18+55=73', '32-26=6'),\n    @('80-77=3', '7+74=81'),\n    @('87+5=92', '26+5=31'),\n    @('4+59=63', '80-7=73'),\n    @('60-53=7', '16+29=45'),\n    @('65-37=28', '45+29=74'),\n    @('54-35=19', '38+9=47'),\n    @('2+39=41', '26+59=85'),\n    @('81-78=3', '30-27=3'),\n    @('91-65=26', '46+16=62'),\n    @('17+68=85', '27+39=66'),\n    @('13+69=82', '74+17=91'),\n    @('33+9=42', '35+59=94'),\n    @('85-17=68', '24+48=72'),\n    @('27+28=55', '31-2=29'),\n    @('76-48=28', '81-72=9'),\n    @('90-83=7', '15+26=41'),\n    @('27+58=85', '24-15=9'),\n    @('85-46=39', '85+8=93'),\n)\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph (the date line above the table) -----------------\n$titlePara = $d.Paragraphs.Item(1)\n$titleOld = $Replacements[0][0]\n$titleNew = $Replacements[0][1]\n$titleText = $titlePara.Range.Text.TrimEnd([char]13, [char]7)\nif ($titleText -ne $titleOld) {\n    throw \"Title paragraph: expected '$titleOld' but found '$titleText'\"\n}\nif ($titleOld -ne $titleNew) {\n    $titlePara.Range.Text = $titleNew\n}\n\n# --- Table cells -------------------------------------------------------\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$idx = 1\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $old = $Replacements[$idx][0]\n        $new = $Replacements[$idx][1]\n        $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($cellText -ne $old) {\n            throw \"Cell ($r,$c): expected '$old' but found '$cellText'\"\n        }\n        if ($old -ne $new) {\n            $cell.Range.Text = $new\n        }\n        $idx += 1\n    }\n}\n\nWrite-Output \"done idx=$idx\"\n"}
